# Auto-generated edit script applying scheduled market-data refresh to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 79167.164
$ws.Range("J3").Value = 79167.164
$ws.Range("L3").Value = 79167.164
$ws.Range("N3").Value = -79395.164

$ws.Range("H12").Value = 116.6
$ws.Range("I12").Value = 116.6
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 116.6
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 53.40000000000001
$ws.Range("N12").ClearContents()

$ws.Range("H95").Value = 33200
$ws.Range("J95").Value = 33200
$ws.Range("L95").Value = 33200
$ws.Range("N95").Value = -38692

$ws.Range("H99").Value = 8792
$ws.Range("J99").Value = 60000
$ws.Range("L99").Value = 180000
$ws.Range("N99").Value = -182996

$ws.Range("H102").Value = 79167.164
$ws.Range("J102").Value = 79167.164
$ws.Range("L102").Value = 79167.164
$ws.Range("N102").Value = -85657.164

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H132").Value = 15605.902
$ws.Range("I132").Value = 2661.3333
$ws.Range("K132").Value = 7983.999899999999
$ws.Range("M132").Value = -5453.999899999999

$ws.Range("H137").Value = 3970.611
$ws.Range("I137").Value = 4599.5405
$ws.Range("K137").Value = 13798.6215
$ws.Range("M137").Value = -11248.6215

$ws.Range("H138").Value = 3431.359
$ws.Range("I138").Value = 1752
$ws.Range("J138").Value = 4599.609
$ws.Range("K138").Value = 5256
$ws.Range("L138").Value = 13798.827
$ws.Range("M138").Value = -116
$ws.Range("N138").Value = -24078.827

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3496.6667
$ws.Range("I61").Value = 3280.1428
$ws.Range("K61").Value = 3280.1428
$ws.Range("M61").Value = -3068.1428

$ws.Range("H74").Value = 2667.0264
$ws.Range("I74").Value = 2586.8235
$ws.Range("K74").Value = 2586.8235
$ws.Range("M74").Value = -1712.8235

$ws.Range("H77").Value = 2667.0264
$ws.Range("I77").Value = 2586.8235
$ws.Range("K77").Value = 12934.1175
$ws.Range("M77").Value = -8566.1175

$ws.Range("H95").Value = 75000
$ws.Range("J95").Value = 75000
$ws.Range("L95").Value = 75000
$ws.Range("N95").Value = -80492

$ws.Range("H136").Value = 3496.6667
$ws.Range("I136").Value = 3280.1428
$ws.Range("K136").Value = 9840.428400000001
$ws.Range("M136").Value = -7290.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 8507.333000000001
$ws.Range("I33").Value = 8507.333000000001
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 8507.333000000001
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -8171.333000000001
$ws.Range("N33").ClearContents()

$ws.Range("H105").Value = 3454.8262
$ws.Range("I105").Value = 1816.2307
$ws.Range("K105").Value = 1816.2307
$ws.Range("M105").Value = -69.23070000000007

$ws.Range("H141").Value = 52495
$ws.Range("J141").Value = 52495
$ws.Range("L141").Value = 52495
$ws.Range("N141").Value = -62855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 16913.445
$ws.Range("J28").Value = 17777.625
$ws.Range("L28").Value = 17777.625
$ws.Range("N28").Value = -18267.625

$ws.Range("H31").Value = 2247
$ws.Range("I31").Value = 2247
$ws.Range("K31").Value = 2247
$ws.Range("M31").Value = -1952

$ws.Range("H32").Value = 367.5
$ws.Range("I32").Value = 367.5
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 367.5
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -51.5
$ws.Range("N32").ClearContents()

$ws.Range("H34").Value = 2247
$ws.Range("I34").Value = 2247
$ws.Range("K34").Value = 2247
$ws.Range("M34").Value = -2045

$ws.Range("H99").Value = 10753.5
$ws.Range("I99").Value = 3000
$ws.Range("K99").Value = 3000
$ws.Range("M99").Value = -1502

$ws.Range("H126").Value = 10753.5
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H132").Value = 4360
$ws.Range("I132").Value = 3715.3333
$ws.Range("K132").Value = 11145.9999
$ws.Range("M132").Value = -8615.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 167119.25
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2685

$ws.Range("H73").Value = 167119.25
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -1908

$ws.Range("H127").Value = 7490
$ws.Range("J127").Value = 7490
$ws.Range("L127").Value = 22470
$ws.Range("N127").Value = -32390

$ws.Range("H140").Value = 11905712
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H105").Value = 74394.336
$ws.Range("J105").Value = 74394.336
$ws.Range("L105").Value = 74394.336
$ws.Range("N105").Value = -81382.336

$ws.Range("H132").Value = 3228
$ws.Range("I132").Value = 2521.25
$ws.Range("J132").Value = 4076.1
$ws.Range("K132").Value = 7563.75
$ws.Range("L132").Value = 12228.3
$ws.Range("M132").Value = -5033.75
$ws.Range("N132").Value = -17288.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3982.1614
$ws.Range("J136").Value = 5481
$ws.Range("L136").Value = 16443
$ws.Range("N136").Value = -21543

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5136.1
$ws.Range("I81").Value = 4094.4285
$ws.Range("K81").Value = 8188.857
$ws.Range("M81").Value = -7127.857

$ws.Range("H84").Value = 5136.1
$ws.Range("I84").Value = 4094.4285
$ws.Range("K84").Value = 40944.285
$ws.Range("M84").Value = -35640.285

$ws.Range("H100").Value = 402.25
$ws.Range("I100").Value = 303
$ws.Range("J100").Value = 700
$ws.Range("K100").Value = 606
$ws.Range("L100").Value = 1400
$ws.Range("M100").Value = -65
$ws.Range("N100").Value = -2482

$ws.Range("H126").Value = 1728.3684
$ws.Range("I126").Value = 1333.6875
$ws.Range("K126").Value = 4001.0625
$ws.Range("M126").Value = -1531.0625

$ws.Range("H132").Value = 1892.8667
$ws.Range("I132").Value = 1695.963
$ws.Range("J132").Value = 3665
$ws.Range("K132").Value = 5087.889
$ws.Range("L132").Value = 10995
$ws.Range("M132").Value = -2557.889
$ws.Range("N132").Value = -16055
